$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "4264-MS-EI-DB-DL-REC-RNI-INT-FFConMONTHLYonDAY25-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1-ON-1st"
$newShortName = "426c"

# Update product name (B1) on both sheets - they share the same text.
$ws1.Range("B1").Value = $newProductName
$ws2.Range("B1").Value = $newProductName

# Update short name (B2) on the input sheet - was numeric, now textual.
$ws1.Range("B2").Value = $newShortName

# Move the active cell / selection on the input sheet from B17 to B2.
$ws1.Range("B2").Select()

# Make the output sheet the active tab instead of the input sheet.
$ws2.Select()
